# Applies the Leve profit-recalculation updates described in the commit diff.
# Each block updates the currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) for one Leve row on the given sheet, matched by sheet name + row number.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 4675
$ws.Range("J76").Value = 4500
$ws.Range("L76").Value = 4500
$ws.Range("N76").Value = -5130

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 4675
$ws.Range("J79").Value = 4500
$ws.Range("L79").Value = 4500
$ws.Range("N79").Value = -6684

# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 906.25
$ws.Range("I80").Value = 906.25
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2718.75
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1720.75
$ws.Range("N80").ClearContents()

# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 906.25
$ws.Range("I83").Value = 906.25
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 8156.25
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -3164.25
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 2730.7778
$ws.Range("I63").Value = 1726
$ws.Range("J63").Value = 3233.1667
$ws.Range("K63").Value = 1726
$ws.Range("L63").Value = 3233.1667
$ws.Range("M63").Value = -1040
$ws.Range("N63").Value = -4605.1667

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2730.7778
$ws.Range("I66").Value = 1726
$ws.Range("J66").Value = 3233.1667
$ws.Range("K66").Value = 8630
$ws.Range("L66").Value = 16165.8335
$ws.Range("M66").Value = -5198
$ws.Range("N66").Value = -23029.8335

# Row 97: Ore for Me
$ws.Range("H97").Value = 692.1
$ws.Range("I97").Value = 615.25
$ws.Range("J97").Value = 999.5
$ws.Range("K97").Value = 615.25
$ws.Range("L97").Value = 999.5
$ws.Range("M97").Value = -119.25
$ws.Range("N97").Value = -1991.5

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1949.5
$ws.Range("I102").Value = 1949.5
$ws.Range("K102").Value = 1949.5
$ws.Range("M102").Value = -327.5

# Row 112: Wrapped Knuckles
$ws.Range("H112").Value = 22000
$ws.Range("J112").Value = 22000
$ws.Range("L112").Value = 22000
$ws.Range("N112").Value = -24954

# Row 135: Forgiveness for My Shins
$ws.Range("H135").Value = 46000
$ws.Range("J135").Value = 46000
$ws.Range("L135").Value = 46000
$ws.Range("N135").Value = -56140

$ws = $wb.Worksheets.Item("BSM")
# Row 35: Lancers' Creed
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 94: High Steal
$ws.Range("H94").Value = 2099.9167
$ws.Range("I94").Value = 2601.875
$ws.Range("J94").Value = 1096
$ws.Range("K94").Value = 2601.875
$ws.Range("L94").Value = 1096
$ws.Range("M94").Value = -2150.875
$ws.Range("N94").Value = -1998

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

# Row 122: To Delight a Dancer
$ws.Range("H122").Value = 50000
$ws.Range("I122").Value = 50000
$ws.Range("K122").Value = 50000
$ws.Range("M122").Value = -45100

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 4794.5293
$ws.Range("I134").Value = 4039
$ws.Range("K134").Value = 12117
$ws.Range("M134").Value = -9582

$ws = $wb.Worksheets.Item("CRP")
# Row 25: Bowing to Necessity
$ws.Range("H25").Value = 3500.25
$ws.Range("I25").Value = 3000.3333
$ws.Range("K25").Value = 3000.3333
$ws.Range("M25").Value = -2826.3333

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 2342.9285
$ws.Range("I58").Value = 2300.3635
$ws.Range("J58").Value = 2499
$ws.Range("K58").Value = 2300.3635
$ws.Range("L58").Value = 2499
$ws.Range("M58").Value = -2097.3635
$ws.Range("N58").Value = -2905

# Row 136: Turali Quality
$ws.Range("H136").Value = 2342.9285
$ws.Range("I136").Value = 2300.3635
$ws.Range("J136").Value = 2499
$ws.Range("K136").Value = 6901.0905
$ws.Range("L136").Value = 7497
$ws.Range("M136").Value = -4351.0905
$ws.Range("N136").Value = -12597

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Range("H12").Value = 295.84616
$ws.Range("J12").Value = 379.7143
$ws.Range("L12").Value = 1139.1429
$ws.Range("N12").Value = -1485.1429

# Row 68: Such a Butter Face
$ws.Range("H68").Value = 1116.5
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 924.75
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 2774.25
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -4396.25

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 1116.5
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 924.75
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 8322.75
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -16434.75

# Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 4251.5
$ws.Range("J80").Value = 4251.5
$ws.Range("L80").Value = 12754.5
$ws.Range("N80").Value = -14626.5

# Row 81: It Goes Down Smoothly
$ws.Range("H81").Value = 11507.375
$ws.Range("J81").Value = 12865.571
$ws.Range("L81").Value = 38596.713
$ws.Range("N81").Value = -40842.713

# Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 4251.5
$ws.Range("J83").Value = 4251.5
$ws.Range("L83").Value = 38263.5
$ws.Range("N83").Value = -47623.5

# Row 84: Quenching the Flame (L)
$ws.Range("H84").Value = 11507.375
$ws.Range("J84").Value = 12865.571
$ws.Range("L84").Value = 115790.139
$ws.Range("N84").Value = -127022.139

# Row 87: Soup That Eats Like a Knight
$ws.Range("H87").Value = 314
$ws.Range("I87").Value = 314
$ws.Range("K87").Value = 942
$ws.Range("M87").Value = 306

# Row 90: Like Ma Used to Make (L)
$ws.Range("H90").Value = 314
$ws.Range("I90").Value = 314
$ws.Range("K90").Value = 2826
$ws.Range("M90").Value = 3414

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 4219.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4219.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4219.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -6215.5

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 4219.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4219.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 21097.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -31081.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 2441.5557
$ws.Range("I16").Value = 2860.5715
$ws.Range("J16").Value = 975
$ws.Range("K16").Value = 2860.5715
$ws.Range("L16").Value = 975
$ws.Range("M16").Value = -2690.5715
$ws.Range("N16").Value = -1315

# Row 110: Breeches of Trust
$ws.Range("H110").Value = 41500
$ws.Range("J110").Value = 41500
$ws.Range("L110").Value = 41500
$ws.Range("N110").Value = -49680

# Row 122: Hell on Leather
$ws.Range("H122").Value = 3207.75
$ws.Range("I122").Value = 3207.75
$ws.Range("K122").Value = 9623.25
$ws.Range("M122").Value = -7173.25

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1079
$ws.Range("I132").Value = 973.75
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 2921.25
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -391.25
$ws.Range("N132").Value = -9560
